$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Camp"
$ws.Cells.Item(2,3).Value = "Fpr2"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.1331866666666667
$ws.Cells.Item(2,8).Value = 0.39956
$ws.Cells.Item(2,9).Value = 0.281772779265766
$ws.Cells.Item(2,10).Value = 0.281772779265766
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.559584
$ws.Cells.Item(2,14).Value = 1.678752
$ws.Cells.Item(2,15).Value = 0.003429531108364747
$ws.Cells.Item(2,16).Value = 0.003429531108364747
$ws.Cells.Item(2,17).Value = 0.07452912768
$ws.Cells.Item(2,18).Value = 0.6707621491200001
$ws.Cells.Item(2,19).Value = 0.0009663485119823377
$ws.Cells.Item(2,20).Value = 0.0009663485119823377
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Camp"
$ws.Cells.Item(3,3).Value = "Fpr2"
$ws.Cells.Item(3,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.1331866666666667
$ws.Cells.Item(3,8).Value = 0.39956
$ws.Cells.Item(3,9).Value = 0.281772779265766
$ws.Cells.Item(3,10).Value = 0.281772779265766
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 28.05830633333333
$ws.Cells.Item(3,14).Value = 84.17491899999999
$ws.Cells.Item(3,15).Value = 0.1719613756258118
$ws.Cells.Item(3,16).Value = 0.1719613756258118
$ws.Cells.Item(3,17).Value = 3.736992292848889
$ws.Cells.Item(3,18).Value = 33.63293063564
$ws.Cells.Item(3,19).Value = 0.04845403473644935
$ws.Cells.Item(3,20).Value = 0.04845403473644934
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Camp"
$ws.Cells.Item(4,3).Value = "Fpr2"
$ws.Cells.Item(4,4).Value = "Neutrophils"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.1331866666666667
$ws.Cells.Item(4,8).Value = 0.39956
$ws.Cells.Item(4,9).Value = 0.281772779265766
$ws.Cells.Item(4,10).Value = 0.281772779265766
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 131.297198
$ws.Cells.Item(4,14).Value = 393.891594
$ws.Cells.Item(4,15).Value = 0.8046831663916869
$ws.Cells.Item(4,16).Value = 0.8046831663916869
$ws.Cells.Item(4,17).Value = 17.48703614429333
$ws.Cells.Item(4,18).Value = 157.38332529864
$ws.Cells.Item(4,19).Value = 0.2267378122225625
$ws.Cells.Item(4,20).Value = 0.2267378122225625
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Camp"
$ws.Cells.Item(5,3).Value = "Fpr2"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.1331866666666667
$ws.Cells.Item(5,8).Value = 0.39956
$ws.Cells.Item(5,9).Value = 0.281772779265766
$ws.Cells.Item(5,10).Value = 0.281772779265766
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 3.251240333333333
$ws.Cells.Item(5,14).Value = 9.753720999999999
$ws.Cells.Item(5,15).Value = 0.01992592687413657
$ws.Cells.Item(5,16).Value = 0.01992592687413657
$ws.Cells.Item(5,17).Value = 0.4330218625288889
$ws.Cells.Item(5,18).Value = 3.89719676276
$ws.Cells.Item(5,19).Value = 0.005614583794771878
$ws.Cells.Item(5,20).Value = 0.005614583794771877
$ws.Cells.Item(6,1).Value = "Neutrophils"
$ws.Cells.Item(6,2).Value = "Camp"
$ws.Cells.Item(6,3).Value = "Fpr2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.1174416666666667
$ws.Cells.Item(6,8).Value = 0.352325
$ws.Cells.Item(6,9).Value = 0.2484622946611547
$ws.Cells.Item(6,10).Value = 0.2484622946611547
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.559584
$ws.Cells.Item(6,14).Value = 1.678752
$ws.Cells.Item(6,15).Value = 0.003429531108364747
$ws.Cells.Item(6,16).Value = 0.003429531108364747
$ws.Cells.Item(6,17).Value = 0.06571847759999999
$ws.Cells.Item(6,18).Value = 0.5914662984
$ws.Cells.Item(6,19).Value = 0.0008521091687961184
$ws.Cells.Item(6,20).Value = 0.0008521091687961184
$ws.Cells.Item(7,1).Value = "Neutrophils"
$ws.Cells.Item(7,2).Value = "Camp"
$ws.Cells.Item(7,3).Value = "Fpr2"
$ws.Cells.Item(7,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.1174416666666667
$ws.Cells.Item(7,8).Value = 0.352325
$ws.Cells.Item(7,9).Value = 0.2484622946611547
$ws.Cells.Item(7,10).Value = 0.2484622946611547
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 28.05830633333333
$ws.Cells.Item(7,14).Value = 84.17491899999999
$ws.Cells.Item(7,15).Value = 0.1719613756258118
$ws.Cells.Item(7,16).Value = 0.1719613756258118
$ws.Cells.Item(7,17).Value = 3.295214259630555
$ws.Cells.Item(7,18).Value = 29.65692833667499
$ws.Cells.Item(7,19).Value = 0.04272591798107797
$ws.Cells.Item(7,20).Value = 0.04272591798107796
$ws.Cells.Item(8,1).Value = "Neutrophils"
$ws.Cells.Item(8,2).Value = "Camp"
$ws.Cells.Item(8,3).Value = "Fpr2"
$ws.Cells.Item(8,4).Value = "Neutrophils"
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 0.3333333333333333
$ws.Cells.Item(8,7).Value = 0.1174416666666667
$ws.Cells.Item(8,8).Value = 0.352325
$ws.Cells.Item(8,9).Value = 0.2484622946611547
$ws.Cells.Item(8,10).Value = 0.2484622946611547
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 131.297198
$ws.Cells.Item(8,14).Value = 393.891594
$ws.Cells.Item(8,15).Value = 0.8046831663916869
$ws.Cells.Item(8,16).Value = 0.8046831663916869
$ws.Cells.Item(8,17).Value = 15.41976176178333
$ws.Cells.Item(8,18).Value = 138.77785585605
$ws.Cells.Item(8,19).Value = 0.1999334259968823
$ws.Cells.Item(8,20).Value = 0.1999334259968823
$ws.Cells.Item(9,1).Value = "Neutrophils"
$ws.Cells.Item(9,2).Value = "Camp"
$ws.Cells.Item(9,3).Value = "Fpr2"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 0.3333333333333333
$ws.Cells.Item(9,7).Value = 0.1174416666666667
$ws.Cells.Item(9,8).Value = 0.352325
$ws.Cells.Item(9,9).Value = 0.2484622946611547
$ws.Cells.Item(9,10).Value = 0.2484622946611547
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.251240333333333
$ws.Cells.Item(9,14).Value = 9.753720999999999
$ws.Cells.Item(9,15).Value = 0.01992592687413657
$ws.Cells.Item(9,16).Value = 0.01992592687413657
$ws.Cells.Item(9,17).Value = 0.3818310834805555
$ws.Cells.Item(9,18).Value = 3.436479751324999
$ws.Cells.Item(9,19).Value = 0.004950841514398343
$ws.Cells.Item(9,20).Value = 0.004950841514398342
$ws.Cells.Item(10,1).Value = "Resolving-Mac"
$ws.Cells.Item(10,2).Value = "Camp"
$ws.Cells.Item(10,3).Value = "Fpr2"
$ws.Cells.Item(10,4).Value = "FAPs"
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(10,7).Value = 0.2220456666666667
$ws.Cells.Item(10,8).Value = 0.666137
$ws.Cells.Item(10,9).Value = 0.4697649260730792
$ws.Cells.Item(10,10).Value = 0.4697649260730792
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.559584
$ws.Cells.Item(10,14).Value = 1.678752
$ws.Cells.Item(10,15).Value = 0.003429531108364747
$ws.Cells.Item(10,16).Value = 0.003429531108364747
$ws.Cells.Item(10,17).Value = 0.124253202336
$ws.Cells.Item(10,18).Value = 1.118278821024
$ws.Cells.Item(10,19).Value = 0.001611073427586291
$ws.Cells.Item(10,20).Value = 0.001611073427586291
$ws.Cells.Item(11,1).Value = "Resolving-Mac"
$ws.Cells.Item(11,2).Value = "Camp"
$ws.Cells.Item(11,3).Value = "Fpr2"
$ws.Cells.Item(11,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.2220456666666667
$ws.Cells.Item(11,8).Value = 0.666137
$ws.Cells.Item(11,9).Value = 0.4697649260730792
$ws.Cells.Item(11,10).Value = 0.4697649260730792
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 28.05830633333333
$ws.Cells.Item(11,14).Value = 84.17491899999999
$ws.Cells.Item(11,15).Value = 0.1719613756258118
$ws.Cells.Item(11,16).Value = 0.1719613756258118
$ws.Cells.Item(11,17).Value = 6.230225335322555
$ws.Cells.Item(11,18).Value = 56.07202801790299
$ws.Cells.Item(11,19).Value = 0.0807814229082845
$ws.Cells.Item(11,20).Value = 0.08078142290828448
$ws.Cells.Item(12,1).Value = "Resolving-Mac"
$ws.Cells.Item(12,2).Value = "Camp"
$ws.Cells.Item(12,3).Value = "Fpr2"
$ws.Cells.Item(12,4).Value = "Neutrophils"
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 0.3333333333333333
$ws.Cells.Item(12,7).Value = 0.2220456666666667
$ws.Cells.Item(12,8).Value = 0.666137
$ws.Cells.Item(12,9).Value = 0.4697649260730792
$ws.Cells.Item(12,10).Value = 0.4697649260730792
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 131.297198
$ws.Cells.Item(12,14).Value = 393.891594
$ws.Cells.Item(12,15).Value = 0.8046831663916869
$ws.Cells.Item(12,16).Value = 0.8046831663916869
$ws.Cells.Item(12,17).Value = 29.15397386137533
$ws.Cells.Item(12,18).Value = 262.385764752378
$ws.Cells.Item(12,19).Value = 0.3780119281722421
$ws.Cells.Item(12,20).Value = 0.3780119281722421
$ws.Cells.Item(13,1).Value = "Resolving-Mac"
$ws.Cells.Item(13,2).Value = "Camp"
$ws.Cells.Item(13,3).Value = "Fpr2"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0.3333333333333333
$ws.Cells.Item(13,7).Value = 0.2220456666666667
$ws.Cells.Item(13,8).Value = 0.666137
$ws.Cells.Item(13,9).Value = 0.4697649260730792
$ws.Cells.Item(13,10).Value = 0.4697649260730792
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 3.251240333333333
$ws.Cells.Item(13,14).Value = 9.753720999999999
$ws.Cells.Item(13,15).Value = 0.01992592687413657
$ws.Cells.Item(13,16).Value = 0.01992592687413657
$ws.Cells.Item(13,17).Value = 0.7219238273085555
$ws.Cells.Item(13,18).Value = 6.497314445776999
$ws.Cells.Item(13,19).Value = 0.009360501564966348
$ws.Cells.Item(13,20).Value = 0.009360501564966345
